$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 481.375
$ws.Range("I2").Value = 308.5
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 308.5
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -195.5
$ws.Range("N2").Value = -1226

$ws.Range("H62").Value = 2203.0454
$ws.Range("I62").Value = 2577.6365
$ws.Range("J62").Value = 1828.4546
$ws.Range("K62").Value = 2577.6365
$ws.Range("L62").Value = 1828.4546
$ws.Range("M62").Value = -1953.6365
$ws.Range("N62").Value = -3076.4546

$ws.Range("H64").Value = 3837.6191
$ws.Range("I64").Value = 4243.2
$ws.Range("J64").Value = 3241.1765
$ws.Range("K64").Value = 4243.2
$ws.Range("L64").Value = 3241.1765
$ws.Range("M64").Value = -3995.2
$ws.Range("N64").Value = -3737.1765

$ws.Range("H65").Value = 2203.0454
$ws.Range("I65").Value = 2577.6365
$ws.Range("J65").Value = 1828.4546
$ws.Range("K65").Value = 12888.1825
$ws.Range("L65").Value = 9142.273000000001
$ws.Range("M65").Value = -9768.182500000001
$ws.Range("N65").Value = -15382.273

$ws.Range("H67").Value = 3837.6191
$ws.Range("I67").Value = 4243.2
$ws.Range("J67").Value = 3241.1765
$ws.Range("K67").Value = 4243.2
$ws.Range("L67").Value = 3241.1765
$ws.Range("M67").Value = -3385.2
$ws.Range("N67").Value = -4957.1765

$ws.Range("H129").Value = 1056.2565
$ws.Range("J129").Value = 1083.6857
$ws.Range("L129").Value = 3251.0571
$ws.Range("N129").Value = -13251.0571

$ws.Range("H132").Value = 1171.4073
$ws.Range("I132").Value = 891.2653
$ws.Range("J132").Value = 3916.8
$ws.Range("K132").Value = 2673.7959
$ws.Range("L132").Value = 11750.4
$ws.Range("M132").Value = -143.7959000000001
$ws.Range("N132").Value = -16810.4

$ws.Range("H138").Value = 2207.7896
$ws.Range("I138").Value = 1587.0385
$ws.Range("J138").Value = 2530.58
$ws.Range("K138").Value = 4761.1155
$ws.Range("L138").Value = 7591.74
$ws.Range("M138").Value = 378.8845000000001
$ws.Range("N138").Value = -17871.74

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 520.6042
$ws.Range("I2").Value = 419.38635
$ws.Range("J2").Value = 1634
$ws.Range("K2").Value = 419.38635
$ws.Range("L2").Value = 1634
$ws.Range("M2").Value = -306.38635
$ws.Range("N2").Value = -1860

$ws.Range("H32").Value = 6405.2124
$ws.Range("I32").Value = 5302.743
$ws.Range("J32").Value = 20002.334
$ws.Range("K32").Value = 5302.743
$ws.Range("L32").Value = 20002.334
$ws.Range("M32").Value = -5015.743
$ws.Range("N32").Value = -20576.334

$ws.Range("H61").Value = 4240.3
$ws.Range("I61").Value = 4868.0625
$ws.Range("J61").Value = 3522.8572
$ws.Range("K61").Value = 4868.0625
$ws.Range("L61").Value = 3522.8572
$ws.Range("M61").Value = -4656.0625
$ws.Range("N61").Value = -3946.8572

$ws.Range("H74").Value = 2092.68
$ws.Range("I74").Value = 1844.9286
$ws.Range("J74").Value = 2408
$ws.Range("K74").Value = 1844.9286
$ws.Range("L74").Value = 2408
$ws.Range("M74").Value = -970.9286
$ws.Range("N74").Value = -4156

$ws.Range("H77").Value = 2092.68
$ws.Range("I77").Value = 1844.9286
$ws.Range("J77").Value = 2408
$ws.Range("K77").Value = 9224.643
$ws.Range("L77").Value = 12040
$ws.Range("M77").Value = -4856.643
$ws.Range("N77").Value = -20776

$ws.Range("H110").Value = 1620.0303
$ws.Range("I110").Value = 1360.0358
$ws.Range("J110").Value = 3076
$ws.Range("K110").Value = 1360.0358
$ws.Range("L110").Value = 3076
$ws.Range("M110").Value = 684.9641999999999
$ws.Range("N110").Value = -7166

$ws.Range("H116").Value = 520.6042
$ws.Range("I116").Value = 419.38635
$ws.Range("J116").Value = 1634
$ws.Range("K116").Value = 419.38635
$ws.Range("L116").Value = 1634
$ws.Range("M116").Value = 1874.61365
$ws.Range("N116").Value = -6222

$ws.Range("H122").Value = 989260.3
$ws.Range("I122").Value = 1427931.6
$ws.Range("J122").Value = 2249.875
$ws.Range("K122").Value = 4283794.800000001
$ws.Range("L122").Value = 6749.625
$ws.Range("M122").Value = -4281344.800000001
$ws.Range("N122").Value = -11649.625

$ws.Range("H136").Value = 4240.3
$ws.Range("I136").Value = 4868.0625
$ws.Range("J136").Value = 3522.8572
$ws.Range("K136").Value = 14604.1875
$ws.Range("L136").Value = 10568.5716
$ws.Range("M136").Value = -12054.1875
$ws.Range("N136").Value = -15668.5716

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 520.6042
$ws.Range("I3").Value = 419.38635
$ws.Range("J3").Value = 1634
$ws.Range("K3").Value = 419.38635
$ws.Range("L3").Value = 1634
$ws.Range("M3").Value = -305.38635
$ws.Range("N3").Value = -1862

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2871.8772
$ws.Range("I31").Value = 2042.5883
$ws.Range("J31").Value = 4097.7827
$ws.Range("K31").Value = 2042.5883
$ws.Range("L31").Value = 4097.7827
$ws.Range("M31").Value = -1747.5883
$ws.Range("N31").Value = -4687.7827

$ws.Range("H34").Value = 2871.8772
$ws.Range("I34").Value = 2042.5883
$ws.Range("J34").Value = 4097.7827
$ws.Range("K34").Value = 2042.5883
$ws.Range("L34").Value = 4097.7827
$ws.Range("M34").Value = -1840.5883
$ws.Range("N34").Value = -4501.7827

$ws.Range("H35").Value = 3225.6667
$ws.Range("I35").Value = 2379
$ws.Range("J35").Value = 4919
$ws.Range("K35").Value = 2379
$ws.Range("L35").Value = 4919
$ws.Range("M35").Value = -2085
$ws.Range("N35").Value = -5507

$ws.Range("H105").Value = 2596.2083
$ws.Range("I105").Value = 2457.95
$ws.Range("J105").Value = 3287.5
$ws.Range("K105").Value = 2457.95
$ws.Range("L105").Value = 3287.5
$ws.Range("M105").Value = -710.9499999999998
$ws.Range("N105").Value = -6781.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 150
$ws.Range("J11").Value = 200
$ws.Range("L11").Value = 600
$ws.Range("N11").Value = -880

$ws.Range("H34").Value = 1136.3636
$ws.Range("J34").Value = 1386
$ws.Range("L34").Value = 4158
$ws.Range("N34").Value = -4326

$ws.Range("H122").Value = 4530.893
$ws.Range("I122").Value = 444.3846
$ws.Range("K122").Value = 3999.4614
$ws.Range("M122").Value = -1549.4614

$ws.Range("H129").Value = 1377.3636
$ws.Range("I129").Value = 2116.6667
$ws.Range("J129").Value = 1100.125
$ws.Range("K129").Value = 6350.000100000001
$ws.Range("L129").Value = 3300.375
$ws.Range("M129").Value = -1350.000100000001
$ws.Range("N129").Value = -13300.375

$ws.Range("H133").Value = 44641.92
$ws.Range("I133").Value = 130523.75
$ws.Range("J133").Value = 6472.222
$ws.Range("K133").Value = 391571.25
$ws.Range("L133").Value = 19416.666
$ws.Range("M133").Value = -386511.25
$ws.Range("N133").Value = -29536.666

$ws.Range("H136").Value = 3591.2715
$ws.Range("I136").Value = 6868.0625
$ws.Range("J136").Value = 2620.3704
$ws.Range("K136").Value = 20604.1875
$ws.Range("L136").Value = 7861.111199999999
$ws.Range("M136").Value = -15504.1875
$ws.Range("N136").Value = -18061.1112

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H70").Value = 5859.647
$ws.Range("I70").Value = 6411.5
$ws.Range("K70").Value = 6411.5
$ws.Range("M70").Value = -6141.5

$ws.Range("H73").Value = 5859.647
$ws.Range("I73").Value = 6411.5
$ws.Range("K73").Value = 6411.5
$ws.Range("M73").Value = -5475.5

$ws.Range("H97").Value = 1100.6552
$ws.Range("I97").Value = 1119.9615
$ws.Range("J97").Value = 933.3333
$ws.Range("K97").Value = 1119.9615
$ws.Range("L97").Value = 933.3333
$ws.Range("M97").Value = -623.9614999999999
$ws.Range("N97").Value = -1925.3333

$ws.Range("H126").Value = 7630.0527
$ws.Range("I126").Value = 10357.667
$ws.Range("J126").Value = 2954.1428
$ws.Range("K126").Value = 31073.001
$ws.Range("L126").Value = 8862.428400000001
$ws.Range("M126").Value = -28603.001
$ws.Range("N126").Value = -13802.4284

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49370.19
$ws.Range("J7").Value = 1395.6
$ws.Range("L7").Value = 1395.6
$ws.Range("N7").Value = -1619.6

$ws.Range("H32").Value = 6066.6665
$ws.Range("I32").Value = 6066.6665
$ws.Range("K32").Value = 6066.6665
$ws.Range("M32").Value = -5749.6665

$ws.Range("H40").Value = 23811360
$ws.Range("I40").Value = 30304888
$ws.Range("K40").Value = 30304888
$ws.Range("M40").Value = -30304752

$ws.Range("H108").Value = 30626
$ws.Range("J108").Value = 30626
$ws.Range("L108").Value = 30626
$ws.Range("N108").Value = -38306

$ws.Range("H122").Value = 4527957
$ws.Range("I122").Value = 6497566
$ws.Range("J122").Value = 1432857.1
$ws.Range("K122").Value = 19492698
$ws.Range("L122").Value = 4298571.300000001
$ws.Range("M122").Value = -19490248
$ws.Range("N122").Value = -4303471.300000001

$ws.Range("H126").Value = 49370.19
$ws.Range("J126").Value = 1395.6
$ws.Range("L126").Value = 4186.799999999999
$ws.Range("N126").Value = -9126.799999999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 44000
$ws.Range("J125").Value = 44000
$ws.Range("L125").Value = 44000
$ws.Range("N125").Value = -53840

$ws.Range("H126").Value = 815.4
$ws.Range("I126").Value = 573.4375
$ws.Range("K126").Value = 1720.3125
$ws.Range("M126").Value = 749.6875
